# Fruta / hortaliza, semanal
# Inserts two new weekly price rows for "Apio" (Terminal Hortofrutícola
# Agro Chillán) at the top of the existing data block (original row 224),
# pushing the previously existing rows 224-250 down to 226-252.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 224 onward down by inserting two new blank rows at 224.
$ws.Rows.Item(224).Insert()
$ws.Rows.Item(224).Insert()

# --- New row 224 -----------------------------------------------------
$ws.Cells.Item(224, 1).Value = 7
$ws.Cells.Item(224, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(224, 3).Value = "Ñuble"
$ws.Cells.Item(224, 4).Value = 44858
$ws.Cells.Item(224, 5).Value = 16
$ws.Cells.Item(224, 6).Value = 100112017
$ws.Cells.Item(224, 7).Value = "Apio"
$ws.Cells.Item(224, 8).Value = "Americana (o)"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 120
$ws.Cells.Item(224, 11).Value = 8000
$ws.Cells.Item(224, 12).Value = 8500
$ws.Cells.Item(224, 13).Value = 8250
$ws.Cells.Item(224, 14).Value = "`$/docena de matas"
$ws.Cells.Item(224, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(224, 16).Value = 1375
$ws.Cells.Item(224, 17).Value = 6
$ws.Cells.Item(224, 18).Value = "Hortaliza"

# --- New row 225 -----------------------------------------------------
$ws.Cells.Item(225, 1).Value = 7
$ws.Cells.Item(225, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(225, 3).Value = "Ñuble"
$ws.Cells.Item(225, 4).Value = 44858
$ws.Cells.Item(225, 5).Value = 16
$ws.Cells.Item(225, 6).Value = 100112017
$ws.Cells.Item(225, 7).Value = "Apio"
$ws.Cells.Item(225, 8).Value = "Americana (o)"
$ws.Cells.Item(225, 9).Value = "Segunda"
$ws.Cells.Item(225, 10).Value = 80
$ws.Cells.Item(225, 11).Value = 7000
$ws.Cells.Item(225, 12).Value = 7000
$ws.Cells.Item(225, 13).Value = 7000
$ws.Cells.Item(225, 14).Value = "`$/docena de matas"
$ws.Cells.Item(225, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(225, 16).Value = 1167
$ws.Cells.Item(225, 17).Value = 6
$ws.Cells.Item(225, 18).Value = "Hortaliza"

# Make sure the date cells keep the same date/time display format used
# throughout column D.
$ws.Range("D224:D225").NumberFormat = $ws.Range("D223").NumberFormat
